$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 55000
$ws.Range("I63").Value = 55000
$ws.Range("K63").Value = 55000
$ws.Range("M63").Value = -54376

$ws.Range("H66").Value = 55000
$ws.Range("I66").Value = 55000
$ws.Range("K66").Value = 165000
$ws.Range("M66").Value = -161880

$ws.Range("H94").Value = 4766.6665
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H125").Value = 5585
$ws.Range("J125").Value = 5937.222
$ws.Range("L125").Value = 53434.998
$ws.Range("N125").Value = -58354.998

$ws.Range("H138").Value = 6413426.5
$ws.Range("I138").Value = 2479.75
$ws.Range("J138").Value = 7579053
$ws.Range("K138").Value = 7439.25
$ws.Range("L138").Value = 22737159
$ws.Range("M138").Value = -2299.25
$ws.Range("N138").Value = -22747439

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15160035
$ws.Range("I32").Value = 19234936
$ws.Range("J32").Value = 24690.857
$ws.Range("K32").Value = 19234936
$ws.Range("L32").Value = 24690.857
$ws.Range("M32").Value = -19234649
$ws.Range("N32").Value = -25264.857

$ws.Range("H61").Value = 14928177
$ws.Range("I61").Value = 18183976
$ws.Range("J61").Value = 5761.6665
$ws.Range("K61").Value = 18183976
$ws.Range("L61").Value = 5761.6665
$ws.Range("M61").Value = -18183764
$ws.Range("N61").Value = -6185.6665

$ws.Range("H107").Value = 31999.5
$ws.Range("I107").Value = 38999
$ws.Range("J107").Value = 25000
$ws.Range("K107").Value = 38999
$ws.Range("L107").Value = 25000
$ws.Range("M107").Value = -35159
$ws.Range("N107").Value = -32680

$ws.Range("H122").Value = 2889.524
$ws.Range("I122").Value = 1906.4445
$ws.Range("K122").Value = 5719.333500000001
$ws.Range("M122").Value = -3269.333500000001

$ws.Range("H132").Value = 22223806
$ws.Range("I132").Value = 1604.1628
$ws.Range("K132").Value = 4812.4884
$ws.Range("M132").Value = -2282.4884

$ws.Range("H136").Value = 14928177
$ws.Range("I136").Value = 18183976
$ws.Range("J136").Value = 5761.6665
$ws.Range("K136").Value = 54551928
$ws.Range("L136").Value = 17284.9995
$ws.Range("M136").Value = -54549378
$ws.Range("N136").Value = -22384.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2192.3704
$ws.Range("I20").Value = 2553
$ws.Range("J20").Value = 1471.1111
$ws.Range("K20").Value = 2553
$ws.Range("L20").Value = 1471.1111
$ws.Range("M20").Value = -2306
$ws.Range("N20").Value = -1965.1111

$ws.Range("H58").Value = 199999
$ws.Range("J58").Value = 199999
$ws.Range("L58").Value = 199999
$ws.Range("N58").Value = -200587

$ws.Range("H60").Value = 105262
$ws.Range("J60").Value = 105262
$ws.Range("L60").Value = 105262
$ws.Range("N60").Value = -106460

$ws.Range("H80").Value = 4218.6875
$ws.Range("I80").Value = 1861.875
$ws.Range("K80").Value = 1861.875
$ws.Range("M80").Value = -863.875

$ws.Range("H83").Value = 4218.6875
$ws.Range("I83").Value = 1861.875
$ws.Range("K83").Value = 9309.375
$ws.Range("M83").Value = -4317.375

$ws.Range("H94").Value = 1513.1428
$ws.Range("I94").Value = 1513.1428
$ws.Range("K94").Value = 1513.1428
$ws.Range("M94").Value = -1062.1428

$ws.Range("H105").Value = 1884.7858
$ws.Range("I105").Value = 1850
$ws.Range("K105").Value = 1850
$ws.Range("M105").Value = -103

$ws.Range("H107").Value = 5947.9165
$ws.Range("I107").Value = 4174
$ws.Range("K107").Value = 4174
$ws.Range("M107").Value = -2254

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2838.7
$ws.Range("I58").Value = 1565
$ws.Range("J58").Value = 4749.25
$ws.Range("K58").Value = 1565
$ws.Range("L58").Value = 4749.25
$ws.Range("M58").Value = -1362
$ws.Range("N58").Value = -5155.25

$ws.Range("H74").Value = 42142.5

$ws.Range("H77").Value = 42142.5

$ws.Range("H105").Value = 10199.571
$ws.Range("I105").Value = 3031.5557
$ws.Range("J105").Value = 23102
$ws.Range("K105").Value = 3031.5557
$ws.Range("L105").Value = 23102
$ws.Range("M105").Value = -1284.5557
$ws.Range("N105").Value = -26596

$ws.Range("H120").Value = 38213.273
$ws.Range("J120").Value = 37192.715
$ws.Range("L120").Value = 37192.715
$ws.Range("N120").Value = -44450.715

$ws.Range("H132").Value = 2783.25
$ws.Range("I132").Value = 2145.842
$ws.Range("J132").Value = 5205.4
$ws.Range("K132").Value = 6437.526
$ws.Range("L132").Value = 15616.2
$ws.Range("M132").Value = -3907.526
$ws.Range("N132").Value = -20676.2

$ws.Range("H134").Value = 1384.1177
$ws.Range("I134").Value = 1283.1875
$ws.Range("K134").Value = 3849.5625
$ws.Range("M134").Value = -1314.5625

$ws.Range("H136").Value = 2838.7
$ws.Range("I136").Value = 1565
$ws.Range("J136").Value = 4749.25
$ws.Range("K136").Value = 4695
$ws.Range("L136").Value = 14247.75
$ws.Range("M136").Value = -2145
$ws.Range("N136").Value = -19347.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 29724.75
$ws.Range("J111").Value = 16633
$ws.Range("L111").Value = 16633
$ws.Range("N111").Value = -22767

$ws.Range("H126").Value = 13049070
$ws.Range("I126").Value = 8340096.5
$ws.Range("K126").Value = 25020289.5
$ws.Range("M126").Value = -25017819.5

$ws.Range("H132").Value = 2868.4211
$ws.Range("I132").Value = 2324.1765
$ws.Range("K132").Value = 6972.529500000001
$ws.Range("M132").Value = -4442.529500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 489.4762
$ws.Range("I16").Value = 414.73685
$ws.Range("J16").Value = 1199.5
$ws.Range("K16").Value = 414.73685
$ws.Range("L16").Value = 1199.5
$ws.Range("M16").Value = -244.73685
$ws.Range("N16").Value = -1539.5

$ws.Range("H22").Value = 3606.818
$ws.Range("I22").Value = 1475.6
$ws.Range("K22").Value = 1475.6
$ws.Range("M22").Value = -1180.6

$ws.Range("H27").Value = 3606.818
$ws.Range("I27").Value = 1475.6
$ws.Range("K27").Value = 1475.6
$ws.Range("M27").Value = -1368.6

$ws.Range("H46").Value = 2326.5278
$ws.Range("I46").Value = 1063.1578
$ws.Range("J46").Value = 3738.5293
$ws.Range("K46").Value = 1063.1578
$ws.Range("L46").Value = 3738.5293
$ws.Range("M46").Value = -875.1578
$ws.Range("N46").Value = -4114.5293

$ws.Range("H55").Value = 664.96
$ws.Range("I55").Value = 352.94116
$ws.Range("J55").Value = 1328
$ws.Range("K55").Value = 352.94116
$ws.Range("L55").Value = 1328
$ws.Range("M55").Value = -179.94116
$ws.Range("N55").Value = -1674

$ws.Range("H127").Value = 68800
$ws.Range("J127").Value = 68800
$ws.Range("L127").Value = 68800
$ws.Range("N127").Value = -78720

$ws.Range("H132").Value = 153847600
$ws.Range("I132").Value = 1712.8182
$ws.Range("K132").Value = 5138.4546
$ws.Range("M132").Value = -2608.4546

$ws.Range("H136").Value = 4025.348
$ws.Range("I136").Value = 4025.348
$ws.Range("K136").Value = 12076.044
$ws.Range("M136").Value = -9526.044

$ws.Range("H137").Value = 23722.5
$ws.Range("J137").Value = 27250
$ws.Range("L137").Value = 27250
$ws.Range("N137").Value = -37450

$ws.Range("H139").Value = 50833.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 70618.5
$ws.Range("J103").Value = 70618.5
$ws.Range("L103").Value = 70618.5
$ws.Range("N103").Value = -72962.5

$ws.Range("H107").Value = 386.46667
$ws.Range("I107").Value = 324.66666
$ws.Range("J107").Value = 479.16666
$ws.Range("K107").Value = 973.9999799999999
$ws.Range("L107").Value = 1437.49998
$ws.Range("M107").Value = 946.0000200000001
$ws.Range("N107").Value = -5277.499980000001

$ws.Range("H132").Value = 3690.0725
$ws.Range("I132").Value = 4160.278
$ws.Range("K132").Value = 12480.834
$ws.Range("M132").Value = -9950.834000000001

$ws.Range("H136").Value = 1895.9459
$ws.Range("I136").Value = 1981.2667
$ws.Range("K136").Value = 5943.800099999999
$ws.Range("M136").Value = -3393.800099999999
